$wb = $excel.ActiveWorkbook

# Add a new worksheet positioned after the last existing sheet (mirrors
# Excel's Worksheets.Add(Before, After) signature) and give it the name
# used by the new "OpenAccountTest" data-driven test case.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "OpenAccountTest"

# Populate the header row and the single data row used by the test.
$ws.Range("A1").Value = "customer"
$ws.Range("B1").Value = "currency"
$ws.Range("A2").Value = "Katya Smith"
$ws.Range("B2").Value = "Dollar"

# Make the new sheet the active tab / selected cell, matching the
# workbook view state captured for the new sheet.
[void]$ws.Range("E7").Select()
